# Hyperparameter_Search.xlsx edit
# - no shuffle, scale wind speeds as well, arrange benchmarks

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---- Header row tweaks (rename "HuberLoss 1st fold" -> "..., test 1"; add test 2 column) ----
$ws1.Range("G1").Value = "HuberLoss 1st fold, test 1"
$ws1.Range("H1").Value = "HuberLoss 1st fold, test 2"

# ---- New H-column values for existing rows ----
$ws1.Range("H19").Value = 18.55
$ws1.Range("H22").Value = 18.43

# ---- Fix typos in configuration descriptions (Normalisaiton -> Normalisaton) ----
$ws1.Range("F20").Value = "4 Linear, 3 Batch Normalisaton, ReLU, Dropout  in between, sizes 256-128-64"
$ws1.Range("F21").Value = "3 Linear, 2 Batch Normalisaton, ReLU, Dropout  in between, sizes 256-64"

# ---- New rows 23 and 24 ----
$ws1.Range("F23").Value = "4 Linear, 3 LeakyReLU in between, sizes 256-128-64"
$ws1.Range("H23").Value = 18.76

$ws1.Range("F24").Value = "5 Linear, 4 LeakyReLU in between, sizes 256-128-64-32"
$ws1.Range("H24").Value = 18.38

# ---- Highlight (yellow fill) the benchmark rows used for scaling ----
$ws1.Range("A8").Interior.Color = 65535
$ws1.Range("B13").Interior.Color = 65535
$ws1.Range("C14").Interior.Color = 65535
$ws1.Range("F19").Interior.Color = 65535

# ---- Column widths on Sheet1 (closest achievable via character-width rounding) ----
$ws1.Columns.Item(1).ColumnWidth = 12
$ws1.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws1.Columns.Item(3).ColumnWidth = 12.833333333333334
$ws1.Columns.Item(4).ColumnWidth = 15.5
$ws1.Columns.Item(5).ColumnWidth = 13.5
$ws1.Columns.Item(6).ColumnWidth = 62
$ws1.Columns.Item(7).ColumnWidth = 17.333333333333332
$ws1.Columns.Item(8).ColumnWidth = 16.5

# ---- Table on Sheet1 ----
$tbl1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:H24"), $null, 1)
$tbl1.Name = "Table2"

# ---- Add Sheet2 after Sheet1 ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Benchmark"
$ws2.Range("B1").Value = "Huber"
$ws2.Range("C1").Value = "MAE / MW"
$ws2.Range("D1").Value = "MSE /MW²"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 27.33
$ws2.Range("C2").Value = 27.79
$ws2.Range("D2").Value = 2986

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 25.53
$ws2.Range("C3").Value = 25.99
$ws2.Range("D3").Value = 2728.96

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 28.01
$ws2.Range("C4").Value = 28.48
$ws2.Range("D4").Value = 3127.33

# ---- Column widths on Sheet2 ----
$ws2.Columns.Item(1).ColumnWidth = 15
$ws2.Columns.Item(3).ColumnWidth = 10
$ws2.Columns.Item(4).ColumnWidth = 10.5

# ---- Table on Sheet2 ----
$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:D4"), $null, 1)
$tbl2.Name = "Table1"

# ---- Selections matching the final author state ----
$ws2.Range("D27").Select()
$ws1.Range("F24").Select()
$ws1.Select()

Write-Host "edit complete"
